$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "asu@gmail.com"
$ws.Range("B3").Value = "ed9c139d8f537c8b631ccc0258c9b570a9fdf723329871a1560c90ba0c95d439"
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = "user"
